# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new labels in AD1:AF1, matching the style of the existing
#     header cells (bold, bordered, centered). Copy the format from an
#     existing header cell (A1) so we reuse the same style index rather than
#     creating a new one.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows 2-51: every player on the roster shares the team's season
#     record of 59 wins, 103 losses, 0 ties.
$wins = 59
$losses = 103
$ties = 0

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
